$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 140
$ws.Cells.Item(140, 1).Value = 45970
$ws.Cells.Item(140, 2).Value = "四方坪站充电量(kw)"
$row140 = @(611.32899999999984, 1426.4450000000004, 595.00599999999997, 573.70699999999999, 319.41499999999996, 593.57400000000007, 333.37800000000004, 375.72999999999996, 184.06399999999999, 181.57499999999999, 192.989, 276.928, 536.27499999999998, 1518.5669999999998, 586.20500000000004, 524.69899999999996, 466.86900000000003, 127.21300000000001, 127.75, 67.597999999999999, 100.16, 56.034999999999997, 105.12, 182.86100000000002)
$col = 3
foreach ($v in $row140) {
    $ws.Cells.Item(140, $col).Value = $v
    $col = $col + 1
}

# Row 141
$ws.Cells.Item(141, 1).Value = 45970
$ws.Cells.Item(141, 2).Value = "高岭站充电量(kw)"
$row141 = @(451.00600000000003, 279.55300000000005, 73.14, 86.156000000000006, 53.258000000000003, 195.351, 90.531000000000006, 85.1, 237.809, 150.55500000000001, 236.876, 241.12899999999999, 351.04700000000003, 554.33299999999997, 368.04200000000003, 278.73500000000001, 213.804, 76.393000000000001, 54.21, 11.93, 11.875, 41.265999999999998, 61.039000000000001, 57.561999999999998)
$col = 3
foreach ($v in $row141) {
    $ws.Cells.Item(141, $col).Value = $v
    $col = $col + 1
}

[void]$ws.Range("F143").Select()